$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify steel description: remove the "/RME" segment from the
# industrial mapping text in cell B2.
$cell = $ws.Range("B2")
$text = $cell.Value2
$newText = $text -replace "/RME/H:1", "/H:1"
$cell.Value = $newText

# Wrap the long descriptive text and size the row to fit it.
$cell.WrapText = $true
$ws.Rows(2).RowHeight = 409.6

# Restore the selection state (A2:A7, active cell A7) as left by the author.
$ws.Range("A7").Activate()
$ws.Range("A2:A7").Select()
